$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "MODEL_CONDITION" header typo/rename before the shift so we
# don't need to re-locate it afterwards.
$ws.Cells.Item(1, 5).Value = "MODELCONDITION"

# Delete column A entirely - shifts B:F left to A:E, matching the diff
# (old col A values of 2/9 are removed, everything else shifts over).
$ws.Columns.Item(1).Delete()
